$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(59.93999862670898, 46.38000106811523, 59.93999862670898, 43.86000061035156, 41948874)
    3  = @(59.93999862670898, 46.38000106811523, 59.93999862670898, 43.86000061035156, 41948874)
    4  = @(59.93999862670898, 46.38000106811523, 59.93999862670898, 43.86000061035156, 41948874)
    5  = @(59.93999862670898, 46.38000106811523, 59.93999862670898, 43.86000061035156, 41948874)
    6  = @(59.93999862670898, 46.38000106811523, 59.93999862670898, 43.86000061035156, 41948874)
    7  = @(59.93999862670898, 46.38000106811523, 59.93999862670898, 43.86000061035156, 41948874)
    8  = @(59.93999862670898, 46.38000106811523, 59.93999862670898, 43.86000061035156, 41948874)
    9  = @(46.5, 42.41999816894531, 47.40000152587891, 35.70000076293945, 41948874)
    10 = @(59.93999862670898, 34.79999923706055, 59.93999862670898, 31.79999923706055, 41948874)
    11 = @(39.77999877929688, 28.20000076293945, 42.54000091552734, 25.86599922180176, 41948874)
    12 = @(15.89999961853027, 20.10000038146973, 20.63999938964844, 15.42000007629394, 41948874)
    13 = @(13.31999969482422, 12.11999988555908, 14.76000022888184, 10.38000011444092, 41948874)
    14 = @(3.420000076293945, 3.599999904632568, 5.165999889373779, 3.059999942779541, 41948874)
    15 = @(2.940000057220459, 2.400000095367432, 3.180000066757202, 2.124000072479248, 41948874)
    16 = @(1.925999999046326, 1.799999952316284, 2.009999990463257, 1.429999947547913, 41948874)
    17 = @(0.8700000047683716, 0.5099999904632568, 0.9509999752044678, 0.3980000019073486, 41948874)
    18 = @(0.746999979019165, 0.5400000214576721, 0.7490000128746033, 0.5049999952316284, 41948874)
    19 = @(1.879999995231628, 1.350000023841858, 3.138999938964844, 1.309999942779541, 41948874)
    20 = @(1.620000004768372, 2, 2.440000057220459, 1.490000009536743, 41948874)
    21 = @(1.769999980926514, 2.049999952316284, 3.5, 1.677000045776367, 41948874)
    22 = @(11.39999961853027, 8.510000228881836, 12.60000038146973, 6.590000152587891, 41948874)
    23 = @(5.989999771118164, 6.079999923706055, 6.800000190734863, 4.829999923706055, 41948874)
    24 = @(8.699999809265137, 8.180000305175781, 13, 7.650000095367432, 41948874)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
    $ws.Cells.Item($row, 6).Value = $vals[2]
    $ws.Cells.Item($row, 7).Value = $vals[3]
    $ws.Cells.Item($row, 8).Value = $vals[4]
    $ws.Cells.Item($row, 9).Value = "PDYN"
}
